$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range('D2').Value = '35.576.99'
$ws.Range('E2').Value = '  +3.20%  '
$ws.Range('D3').Value = '1.843.00'
$ws.Range('E3').Value = '  +2.18%  '
$ws.Range('E4').Value = '  +0.17%  '
$cell = $ws.Range('D5')
$origStyle = $cell.Style
$cell.NumberFormat = '@'
$cell.Value = '231.95'
$cell.Style = $origStyle
$ws.Range('E5').Value = '  +3.29%  '
$cell = $ws.Range('D6')
$origStyle = $cell.Style
$cell.NumberFormat = '@'
$cell.Value = '0.618'
$cell.Style = $origStyle
$ws.Range('E6').Value = '  +2.66%  '
$ws.Range('E7').Value = '  +0.22%  '
$cell = $ws.Range('D8')
$origStyle = $cell.Style
$cell.NumberFormat = '@'
$cell.Value = '43.84'
$cell.Style = $origStyle
$ws.Range('E8').Value = '  +10.67%  '
$cell = $ws.Range('D9')
$origStyle = $cell.Style
$cell.NumberFormat = '@'
$cell.Value = '0.312'
$cell.Style = $origStyle
$ws.Range('E9').Value = '  +8.54%  '
$cell = $ws.Range('D10')
$origStyle = $cell.Style
$cell.NumberFormat = '@'
$cell.Value = '0.0705'
$cell.Style = $origStyle
$ws.Range('E10').Value = '  +5.38%  '
$ws.Range('E11').Value = '  +2.79%  '
$ws.Range('D12').Value = '2.108.03'
$ws.Range('E12').Value = '  +2.01%  '
$ws.Range('B13').Value = 'WrappedEther'
$ws.Range('C13').Value = 'https://coinranking.com/coin/Mtfb0obXVh59u+wrappedether-weth'
$ws.Range('D13').Value = '1.856.04'
$ws.Range('E13').Value = '  +3.10%  '
$ws.Range('B14').Value = 'Chainlink'
$ws.Range('C14').Value = 'https://coinranking.com/coin/VLqpJwogdhHNb+chainlink-link'
$cell = $ws.Range('D14')
$origStyle = $cell.Style
$cell.NumberFormat = '@'
$cell.Value = '11.30'
$cell.Style = $origStyle
$ws.Range('E14').Value = '  +2.99%  '
$ws.Range('B15').Value = 'Polygon'
$ws.Range('C15').Value = 'https://coinranking.com/coin/uW2tk-ILY0ii+polygon-matic'
$cell = $ws.Range('D15')
$origStyle = $cell.Style
$cell.NumberFormat = '@'
$cell.Value = '0.676'
$cell.Style = $origStyle
$ws.Range('E15').Value = '  +7.18%  '
$cell = $ws.Range('D16')
$origStyle = $cell.Style
$cell.NumberFormat = '@'
$cell.Value = '4.75'
$cell.Style = $origStyle
$ws.Range('E16').Value = '  +8.58%  '
$ws.Range('D17').Value = '35.541.69'
$ws.Range('E17').Value = '  +3.15%  '
$cell = $ws.Range('D18')
$origStyle = $cell.Style
$cell.NumberFormat = '@'
$cell.Value = '70.50'
$cell.Style = $origStyle
$ws.Range('E18').Value = '  +3.48%  '
$ws.Range('D19').Value = '0.0₃0804'
$ws.Range('E19').Value = '  +4.73%  '
$cell = $ws.Range('D20')
$origStyle = $cell.Style
$cell.NumberFormat = '@'
$cell.Value = '245.48'
$cell.Style = $origStyle
$ws.Range('E20').Value = '  +2.55%  '
$cell = $ws.Range('D21')
$origStyle = $cell.Style
$cell.NumberFormat = '@'
$cell.Value = '12.07'
$cell.Style = $origStyle
$ws.Range('E21').Value = '  +8.33%  '
$cell = $ws.Range('D22')
$origStyle = $cell.Style
$cell.NumberFormat = '@'
$cell.Value = '4.73'
$cell.Style = $origStyle
$ws.Range('E22').Value = '  +15.80%  '
$ws.Range('E23').Value = '  +0.27%  '
$cell = $ws.Range('D24')
$origStyle = $cell.Style
$cell.NumberFormat = '@'
$cell.Value = '2.22'
$cell.Style = $origStyle
$ws.Range('E24').Value = '  +2.24%  '
$cell = $ws.Range('D25')
$origStyle = $cell.Style
$cell.NumberFormat = '@'
$cell.Value = '172.34'
$cell.Style = $origStyle
$ws.Range('E25').Value = '  +0.44%  '
$cell = $ws.Range('D26')
$origStyle = $cell.Style
$cell.NumberFormat = '@'
$cell.Value = '7.99'
$cell.Style = $origStyle
$ws.Range('E26').Value = '  +4.03%  '
$cell = $ws.Range('D27')
$origStyle = $cell.Style
$cell.NumberFormat = '@'
$cell.Value = '17.87'
$cell.Style = $origStyle
$ws.Range('E27').Value = '  +1.15%  '
$cell = $ws.Range('D28')
$origStyle = $cell.Style
$cell.NumberFormat = '@'
$cell.Value = '0.123'
$cell.Style = $origStyle
$ws.Range('E28').Value = '  +0.84%  '
$cell = $ws.Range('D29')
$origStyle = $cell.Style
$cell.NumberFormat = '@'
$cell.Value = '1.56'
$cell.Style = $origStyle
$ws.Range('E29').Value = '  +27.73%  '
$ws.Range('E30').Value = '  +0.16%  '
$ws.Range('D31').Value = '3.347.68'
$ws.Range('E31').Value = '  +37.78%  '
$ws.Range('E32').Value = '  +7.75%  '
$ws.Range('B33').Value = 'Filecoin'
$ws.Range('C33').Value = 'https://coinranking.com/coin/ymQub4fuB+filecoin-fil'
$cell = $ws.Range('D33')
$origStyle = $cell.Style
$cell.NumberFormat = '@'
$cell.Value = '3.95'
$cell.Style = $origStyle
$ws.Range('E33').Value = '  +5.47%  '
$ws.Range('B34').Value = 'InternetComputer(DFINITY)'
$ws.Range('C34').Value = 'https://coinranking.com/coin/aMNLwaUbY+internetcomputerdfinity-icp'
$cell = $ws.Range('D34')
$origStyle = $cell.Style
$cell.NumberFormat = '@'
$cell.Value = '4.09'
$cell.Style = $origStyle
$ws.Range('E34').Value = '  +6.65%  '
$ws.Range('E35').Value = '  +1.39%  '
$cell = $ws.Range('D36')
$origStyle = $cell.Style
$cell.NumberFormat = '@'
$cell.Value = '96.06'
$cell.Style = $origStyle
$ws.Range('E36').Value = '  +17.37%  '
$ws.Range('B37').Value = 'TrustWalletToken'
$ws.Range('C37').Value = 'https://coinranking.com/coin/Hm3OlynlC+trustwallettoken-twt'
$cell = $ws.Range('D37')
$origStyle = $cell.Style
$cell.NumberFormat = '@'
$cell.Value = '1.15'
$cell.Style = $origStyle
$ws.Range('E37').Value = '  +9.10%  '
$ws.Range('B38').Value = 'ImmutableX'
$ws.Range('C38').Value = 'https://coinranking.com/coin/Z96jIvLU7+immutablex-imx'
$cell = $ws.Range('D38')
$origStyle = $cell.Style
$cell.NumberFormat = '@'
$cell.Value = '0.693'
$cell.Style = $origStyle
$ws.Range('E38').Value = '  +8.07%  '
$ws.Range('D39').Value = '1.351.59'
$ws.Range('E39').Value = '  +3.81%  '
$ws.Range('B40').Value = 'RenderToken'
$ws.Range('C40').Value = 'https://coinranking.com/coin/7C4Mh4xy1yDel+rendertoken-rndr'
$cell = $ws.Range('D40')
$origStyle = $cell.Style
$cell.NumberFormat = '@'
$cell.Value = '2.47'
$cell.Style = $origStyle
$ws.Range('E40').Value = '  +7.07%  '
$ws.Range('B41').Value = 'InjectiveProtocol'
$ws.Range('C41').Value = 'https://coinranking.com/coin/PkY9BmsyW+injectiveprotocol-inj'
$cell = $ws.Range('D41')
$origStyle = $cell.Style
$cell.NumberFormat = '@'
$cell.Value = '15.49'
$cell.Style = $origStyle
$ws.Range('E41').Value = '  +11.31%  '
$ws.Range('E42').Value = '  +4.98%  '
$cell = $ws.Range('D43')
$origStyle = $cell.Style
$cell.NumberFormat = '@'
$cell.Value = '1.02'
$cell.Style = $origStyle
$ws.Range('E43').Value = '  +7.31%  '
$cell = $ws.Range('D44')
$origStyle = $cell.Style
$cell.NumberFormat = '@'
$cell.Value = '1.27'
$cell.Style = $origStyle
$ws.Range('E44').Value = '  +4.55%  '
$ws.Range('E45').Value = '  +0.84%  '
$ws.Range('E46').Value = '  +0.96%  '
$ws.Range('E47').Value = '  +9.91%  '
$cell = $ws.Range('D48')
$origStyle = $cell.Style
$cell.NumberFormat = '@'
$cell.Value = '0.0518'
$cell.Style = $origStyle
$ws.Range('E48').Value = '  +0.59%  '
$ws.Range('D49').Value = '2.011.02'
$ws.Range('E49').Value = '  +2.26%  '
$ws.Range('E50').Value = '  +0.28%  '
$cell = $ws.Range('D51')
$origStyle = $cell.Style
$cell.NumberFormat = '@'
$cell.Value = '103.01'
$cell.Style = $origStyle
$ws.Range('E51').Value = '  +0.90%  '
